$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells (Coin names, Links, Price, Volume%) keep their original
# text representation rather than being auto-converted to numbers/dates by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.742.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.348.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.671"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.83"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1000"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.90"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.76"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.696.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.901"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.347.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.710.88"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.03"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.65"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +20.77%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.07%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.21%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.82"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0277"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.32"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +29.32%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +11.33%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.05%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.11"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.25%  "
